$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.885.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "'1.888.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7686"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'242.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.3126"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'25.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'0.07181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "'0.08539"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").Value = "'0.7644"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'1.910.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'5.360"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "'93.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "'6.146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "'29.774.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "'244.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'0.000007817"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.136.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "'8.010"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'0.1647"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").Value = "'9.409"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'162.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'18.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'2.035"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'1.465"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "'1.533"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "'4.512"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "'4.097"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "'0.7425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").Value = "'2.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").Value = "'0.01952"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "'0.4465"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'1.101.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.53%  "
$ws.Range("D43").Value = "'73.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "'6.070"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "'0.8533"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'102.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "'7.675"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "'1.868"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "'2.034.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
